# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet.
# All existing data rows (2-15) get their timestamp bumped from
# 2025-09-25 18:27:05 to 2025-09-25 18:32:54, as described by the commit
# message "Append: 2025-09-25 18:32 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-25 18:27:05"
$newTimestamp = "2025-09-25 18:32:54"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
